# Add four new statistic-summary sheets (AVERAGE (13), STDEVPA (13), MIN (11),
# MAX (11)) at the end of the workbook, mirroring the existing
# AVERAGE/STDEVPA/MIN/MAX sheet families but computed over the original
# "Sapphire:Sophie" sheet range (the first 15 data sheets of the workbook).

$wb = $excel.ActiveWorkbook

# Each new sheet is built by duplicating the nearest same-kind template sheet
# (so it inherits the right column widths, styles, merged cells, hyperlink,
# etc.) and then overwriting its formulas/name.
$templates = @(
    @{ Template = "AVERAGE (12)"; Name = "AVERAGE (13)"; Func = "AVERAGE" },
    @{ Template = "STDEVPA (12)"; Name = "STDEVPA (13)"; Func = "STDEVPA" },
    @{ Template = "MIN (10)";     Name = "MIN (11)";     Func = "MIN" },
    @{ Template = "MAX (10)";     Name = "MAX (11)";     Func = "MAX" }
)

foreach ($t in $templates) {
    $lastIndex = $wb.Worksheets.Count
    $src = $wb.Worksheets.Item($t.Template)
    $src.Copy($null, $wb.Worksheets.Item($lastIndex))

    $new = $wb.Worksheets.Item($wb.Worksheets.Count)
    $new.Name = $t.Name

    $func = $t.Func
    $new.Range("A1").Formula = "=" + $func + "(Sapphire:Sophie!`$A`$1)"
    $new.Range("A8").Formula = "=" + $func + "(Sapphire:Sophie!`$A`$8)"
    $new.Range("B8").Formula = "=" + $func + "(Sapphire:Sophie!`$B`$8)"
    $new.Range("C8").Formula = "=" + $func + "(Sapphire:Sophie!`$C`$8)"
    $new.Range("D8").Formula = "=" + $func + "(Sapphire:Sophie!`$D`$8)"
    $new.Range("B11").Formula = "=" + $func + "(Sapphire:Sophie!`$B`$11)"
}

# The last added sheet ("MAX (11)") becomes the active / selected tab, matching
# the author's workbook view after the edit.
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
